$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(177).Insert()

$ws.Cells.Item(177, 1).Value = 7
$ws.Cells.Item(177, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(177, 3).Value = "Ñuble"
$ws.Cells.Item(177, 4).Value = 44798
$ws.Cells.Item(177, 5).Value = 16
$ws.Cells.Item(177, 6).Value = "Fruta"
$ws.Cells.Item(177, 7).Value = 100104
$ws.Cells.Item(177, 8).Value = "Frutos de pepita"
$ws.Cells.Item(177, 9).Value = 100104005
$ws.Cells.Item(177, 10).Value = "Pera"
$ws.Cells.Item(177, 11).Value = "Packham's Triumph"
$ws.Cells.Item(177, 12).Value = "Primera"
$ws.Cells.Item(177, 13).Value = 80
$ws.Cells.Item(177, 14).Value = 9500
$ws.Cells.Item(177, 15).Value = 10000
$ws.Cells.Item(177, 16).Value = 9750
$ws.Cells.Item(177, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(177, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(177, 19).Value = 609
$ws.Cells.Item(177, 20).Value = 16

Write-Host "done"
